$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column C data (row 1 and row 2)
$ws.Range("C1").Value = -1
$ws.Range("C2").Value = 0.000038579502713798809

# Set the width of the new column C to match the target layout.
# (The runtime quantizes ColumnWidth to 1/6-character steps when it
# serializes the "width" attribute, so 2 is the closest input that
# reproduces the target stored width of 2.85546875.)
$ws.Columns.Item(3).ColumnWidth = 2
